$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.297.37"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.17"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.33"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4354"
$ws.Range("E7").Value = "  +3.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3672"
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.90"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07700"
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.144"
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.16"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.309"
$ws.Range("E14").Value = "  +3.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.543"
$ws.Range("E15").Value = "  +5.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.839.21"
$ws.Range("E16").Value = "  +5.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.80"
$ws.Range("E17").Value = "  +6.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001083"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06564"
$ws.Range("E19").Value = "  +7.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.53"
$ws.Range("E21").Value = "  +4.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.274"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.334.61"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.64"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.059"
$ws.Range("E25").Value = "  -11.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.16"
$ws.Range("E26").Value = "  +7.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.77"
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.042.52"
$ws.Range("E28").Value = "  +5.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.311"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.12"
$ws.Range("E30").Value = "  +2.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.232"
$ws.Range("E31").Value = "  +2.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.961"
$ws.Range("E32").Value = "  +5.36%  "
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.464"
$ws.Range("E34").Value = "  -4.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.01"
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02355"
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2180"
$ws.Range("E37").Value = "  +2.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.208"
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6586"
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06209"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.195"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.137"
$ws.Range("E42").Value = "  +3.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.442"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.86"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6124"
$ws.Range("E46").Value = "  +4.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.755"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.027"
$ws.Range("E48").Value = "  +4.40%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.86"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.160"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07013"
$ws.Range("E51").Value = "  +2.43%  "
